$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "_old" / "_new" header-suffix columns to "_FV2404" / "_FV2410"
$fv2404Headers = @("Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404","Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404")
$fv2410Headers = @("Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410","Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410")

for ($i = 0; $i -lt 10; $i++) {
    $colLeft = $i + 1          # A..J
    $colRight = $i + 12        # L..U
    $ws.Cells.Item(1, $colLeft).Value = $fv2404Headers[$i]
    $ws.Cells.Item(1, $colRight).Value = $fv2410Headers[$i]
}

# Wrap the used range in an Excel table (adds xl/tables/table1.xml + tableParts)
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U85"), $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (row 1) via a split pane
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
